$d = $word.ActiveDocument

# 1. Insert three new bullet paragraphs immediately before the
#    "Developed and deployed custom analytical tools..." bullet in the
#    Siege Analytics / Advanced Data Analysis and Statistical Modeling section.
$target = "Developed and deployed custom analytical tools and algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*$target*") {
        $insertionPoint = $d.Range($p.Range.Start, $p.Range.Start)
        $newText = "• Developed meta-analytical techniques that identified systematic data quality issues across 20+ years of voter registration data`r" + `
                   "• Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters`r" + `
                   "• Created fraud detection systems analyzing 5+ terabyte datasets, uncovering demographic miscoding patterns across 2,000+ precincts`r"
        $insertionPoint.InsertBefore($newText)
        break
    }
}

# 2. Remove the now-redundant "Created fraud detection systems for campaign
#    finance data analysis across multi-terabyte datasets" bullet.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Created fraud detection systems for campaign finance data analysis across multi-terabyte datasets*") {
        $p.Range.Delete()
        break
    }
}
